# Apply updated dual-variable values produced by the working branch-and-price
# run with L-shaped subproblems.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "u_MAB"
# ---------------------------------------------------------------------------
$wsMAB = $wb.Worksheets.Item("u_MAB")

$wsMAB.Range("A16").Value = 0.05412794201161033
$wsMAB.Range("B16").Value = 0.03586158011550228

$wsMAB.Range("A22").Value = 1.185135361747168

$wsMAB.Range("A23").Value = 0.2640642615044496
$wsMAB.Range("B23").Value = 0.4131639868588319

$wsMAB.Range("A25").Value = 0.2162527072055038
$wsMAB.Range("B25").Value = 0.07835848448175436

$wsMAB.Range("A27").Value = 0.09515418782713028
$wsMAB.Range("B27").Value = 0.1041391036055482

$wsMAB.Range("B41").Value = 0.1644086401790201

$wsMAB.Range("A47").Value = 0.006629605030329028

$wsMAB.Range("A49").Value = 0.1155405216694803
$wsMAB.Range("B49").Value = 0.2502207903436298

$wsMAB.Range("B51").Value = 0.0860122782694944

$wsMAB.Range("A52").Value = 0.05182702263477312

$wsMAB.Range("B60").Value = 0
$wsMAB.Range("B61").Value = 0

# ---------------------------------------------------------------------------
# Sheet "u_EOH"
# ---------------------------------------------------------------------------
$wsEOH = $wb.Worksheets.Item("u_EOH")

$wsEOH.Range("A2").Value = -0.563948174827347
$wsEOH.Range("A3").Value = -0.5047623472507067

# ---------------------------------------------------------------------------
# Sheet "v_l"
# ---------------------------------------------------------------------------
$wsVL = $wb.Worksheets.Item("v_l")

$wsVL.Range("A2").Value = 3935940.064094177
$wsVL.Range("A3").Value = 4594514.000345765
$wsVL.Range("A4").Value = 0
